$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.5368887870543348
$ws.Range("C4").Value = 0.53
$ws.Range("D4").Value = 0.5881833695699248
$ws.Range("E4").Value = 0.5814999999999999
$ws.Range("F4").Value = 0.6654802465277141
$ws.Range("G4").Value = 0.986
$ws.Range("H4").Value = 0.5022520498091527
$ws.Range("I4").Value = 0.5044999999999999
$ws.Range("J4").Value = 0.5837208477215682
$ws.Range("K4").Value = 0.602
$ws.Range("L4").Value = 0.591064732320134
$ws.Range("M4").Value = 0.602
$ws.Range("B5").Value = 0.6775616024809574
$ws.Range("C5").Value = 0.795
$ws.Range("D5").Value = 0.6006652472605893
$ws.Range("E5").Value = 0.6199999999999999
$ws.Range("F5").Value = 0.6581324052970674
$ws.Range("G5").Value = 0.951
$ws.Range("H5").Value = 0.5034250373314075
$ws.Range("I5").Value = 0.507
$ws.Range("J5").Value = 0.622681016559759
$ws.Range("K5").Value = 0.642
$ws.Range("L5").Value = 0.6242898569253035
$ws.Range("M5").Value = 0.6264999999999998
$ws.Range("B6").Value = 0.2315652181000164
$ws.Range("C6").Value = 0.158
$ws.Range("D6").Value = 0.5081048387096774
$ws.Range("E6").Value = 0.5135
$ws.Range("F6").Value = 0.6584475508878966
$ws.Range("G6").Value = 0.9710000000000001
$ws.Range("H6").Value = 0.4982237012481695
$ws.Range("I6").Value = 0.497
$ws.Range("J6").Value = 0.5793743266230057
$ws.Range("K6").Value = 0.58
$ws.Range("L6").Value = 0.5993515714331362
$ws.Range("M6").Value = 0.6180000000000001
$ws.Range("B7").Value = 0.3646551573457102
$ws.Range("C7").Value = 0.375
$ws.Range("D7").Value = 0.3938550618019467
$ws.Range("E7").Value = 0.4065
$ws.Range("F7").Value = 0.1104702787739126
$ws.Range("G7").Value = 0.114
$ws.Range("H7").Value = 0.2728305558383962
$ws.Range("I7").Value = 0.477
$ws.Range("J7").Value = 0.4890675959005309
$ws.Range("K7").Value = 0.509
$ws.Range("L7").Value = 0.5101020151904251
$ws.Range("M7").Value = 0.5130000000000001